$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map (B,C,D,E,F,I,J,K,L,M,N)
$colIndex = @{
    "B" = 2
    "C" = 3
    "D" = 4
    "E" = 5
    "F" = 6
    "I" = 9
    "J" = 10
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
}

# Row data: new vm_pu values for rows 2-25 ("case with 380 kV done")
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.041965429811312; "D" = 1.048149074808328; "E" = 1.045624627833948; "F" = 1.055760267162323; "I" = 1.043773080900968; "J" = 1.047043796655332; "K" = 1.050909726480913; "L" = 1.048392344486278; "M" = 1.058499839548612; "N" = 1.048530718944185 }
    3 = @{ "B" = 1.02; "C" = 1.042919548124774; "D" = 1.048905561855759; "E" = 1.046527558483697; "F" = 1.05670147367993; "I" = 1.044052140462197; "J" = 1.047644136667051; "K" = 1.051478344511386; "L" = 1.049106515277512; "M" = 1.059254225558452; "N" = 1.049131911507581 }
    4 = @{ "B" = 1.02; "C" = 1.043537300010837; "D" = 1.049395354752548; "E" = 1.047112520958926; "F" = 1.057311250256753; "I" = 1.044231557395825; "J" = 1.048032345497218; "K" = 1.051845909373342; "L" = 1.049568704285608; "M" = 1.059742489326609; "N" = 1.049520671638814 }
    5 = @{ "B" = 1.02; "C" = 1.043797091490957; "D" = 1.049601333134822; "E" = 1.047358607010823; "F" = 1.057567779372024; "I" = 1.04430670764455; "J" = 1.048195487482189; "K" = 1.052000344408354; "L" = 1.049763024983213; "M" = 1.059947784420908; "N" = 1.049684045304117 }
    6 = @{ "B" = 1.02; "C" = 1.043840716794135; "D" = 1.049635921848548; "E" = 1.047399935748896; "F" = 1.057610862202591; "I" = 1.0443193094647; "J" = 1.048222876137294; "K" = 1.052026269471847; "L" = 1.049795653218961; "M" = 1.05998225604083; "N" = 1.049711472854254 }
    7 = @{ "B" = 1.02; "C" = 1.04354077101026; "D" = 1.049398106774419; "E" = 1.047115808514828; "F" = 1.057314677309374; "I" = 1.044232562644761; "J" = 1.048034525649179; "K" = 1.051847973293573; "L" = 1.04957130074491; "M" = 1.059745232377054; "N" = 1.049522854886841 }
    8 = @{ "B" = 1.02; "C" = 1.042287800714001; "D" = 1.048404671292275; "E" = 1.04592963058783; "F" = 1.056078195651558; "I" = 1.043867628853703; "J" = 1.047246735898079; "K" = 1.051101969482157; "L" = 1.04863368619694; "M" = 1.058754761195689; "N" = 1.048733946383934 }
    9 = @{ "B" = 1.02; "C" = 1.040082804777141; "D" = 1.046656425710545; "E" = 1.043844892568924; "F" = 1.053905178643152; "I" = 1.043215759964512; "J" = 1.045856662306571; "K" = 1.049784628278622; "L" = 1.0469820898878; "M" = 1.057010432797873; "N" = 1.047341898728482 }
    10 = @{ "B" = 1.02; "C" = 1.038614804116646; "D" = 1.045492556670457; "E" = 1.042458806868994; "F" = 1.05246048705142; "I" = 1.042775289449561; "J" = 1.044928726959065; "K" = 1.0489045788064; "L" = 1.045881486503672; "M" = 1.055848287053488; "N" = 1.046412645606347 }
    11 = @{ "B" = 1.02; "C" = 1.037979626879806; "D" = 1.044988989536832; "E" = 1.041859517728633; "F" = 1.051835880060434; "I" = 1.042583169936641; "J" = 1.044526641545541; "K" = 1.048523086194254; "L" = 1.045405034137484; "M" = 1.055345253332559; "N" = 1.046009989185417 }
    12 = @{ "B" = 1.02; "C" = 1.037743766113302; "D" = 1.044802002918923; "E" = 1.041637050652172; "F" = 1.051604017702999; "I" = 1.042511599408207; "J" = 1.044377247332931; "K" = 1.048381319763738; "L" = 1.045228076872431; "M" = 1.055158432698079; "N" = 1.045860382815889 }
    13 = @{ "B" = 1.02; "C" = 1.037794355763255; "D" = 1.044842109420537; "E" = 1.041684764441476; "F" = 1.051653746408614; "I" = 1.042526960969898; "J" = 1.044409294777193; "K" = 1.048411731986014; "L" = 1.045266033961191; "M" = 1.055198505053419; "N" = 1.045892475771197 }
    14 = @{ "B" = 1.02; "C" = 1.037960129075384; "D" = 1.044973531920947; "E" = 1.041841125758184; "F" = 1.051816711291368; "I" = 1.042577258150866; "J" = 1.044514293421291; "K" = 1.048511369019116; "L" = 1.045390406415015; "M" = 1.055329810093803; "N" = 1.045997623525414 }
    15 = @{ "B" = 1.02; "C" = 1.038062277091532; "D" = 1.045054513730793; "E" = 1.041937483126767; "F" = 1.051917138515214; "I" = 1.042608220252034; "J" = 1.044578981037265; "K" = 1.048572750348393; "L" = 1.045467038848334; "M" = 1.055410715267104; "N" = 1.046062403005222 }
    16 = @{ "B" = 1.02; "C" = 1.038656968843263; "D" = 1.045525985202119; "E" = 1.042498598669511; "F" = 1.052501960359155; "I" = 1.04278801048326; "J" = 1.044955406141194; "K" = 1.048929888322331; "L" = 1.045913109602733; "M" = 1.055881675685652; "N" = 1.046439362675975 }
    17 = @{ "B" = 1.02; "C" = 1.039030131407661; "D" = 1.045821833601215; "E" = 1.042850812057696; "F" = 1.052869059884002; "I" = 1.04290041556908; "J" = 1.045191452365821; "K" = 1.049153798229681; "L" = 1.046192949740736; "M" = 1.056177146417628; "N" = 1.046675744113315 }
    18 = @{ "B" = 1.02; "C" = 1.039247836767938; "D" = 1.045994435081781; "E" = 1.043056338458913; "F" = 1.05308327479536; "I" = 1.042965845129755; "J" = 1.045329106750262; "K" = 1.049284359991966; "L" = 1.046356186905836; "M" = 1.056349507092665; "N" = 1.046813593982772 }
    19 = @{ "B" = 1.02; "C" = 1.039322076441954; "D" = 1.046053294155669; "E" = 1.043126432253044; "F" = 1.053156332108602; "I" = 1.042988132103973; "J" = 1.045376038669616; "K" = 1.049328871182937; "L" = 1.046411848423921; "M" = 1.056408280583304; "N" = 1.046860592550834 }
    20 = @{ "B" = 1.02; "C" = 1.038990089838042; "D" = 1.045790087917038; "E" = 1.042813013960917; "F" = 1.052829664057936; "I" = 1.042888369467825; "J" = 1.045166129666761; "K" = 1.049129779060731; "L" = 1.046162924382737; "M" = 1.056145443377206; "N" = 1.046650385453119 }
    21 = @{ "B" = 1.02; "C" = 1.037911310964024; "D" = 1.044934829586957; "E" = 1.041795077495624; "F" = 1.05176871819723; "I" = 1.042562452644408; "J" = 1.044483375084939; "K" = 1.048482030131207; "L" = 1.045353781313218; "M" = 1.055291143237801; "N" = 1.045966661281478 }
    22 = @{ "B" = 1.02; "C" = 1.037233458767533; "D" = 1.044397446133106; "E" = 1.041155845493028; "F" = 1.05110249597472; "I" = 1.042356328061025; "J" = 1.044053858797271; "K" = 1.048074400067349; "L" = 1.044845147612822; "M" = 1.054754176139298; "N" = 1.045536535031414 }
    23 = @{ "B" = 1.02; "C" = 1.037592761095108; "D" = 1.044682289571256; "E" = 1.041494639686151; "F" = 1.051455593296674; "I" = 1.042465712922913; "J" = 1.04428157612983; "K" = 1.048290526762573; "L" = 1.045114773547825; "M" = 1.055038816559475; "N" = 1.045764575748705 }
    24 = @{ "B" = 1.02; "C" = 1.03900818275782; "D" = 1.045804432307579; "E" = 1.042830093027316; "F" = 1.052847465051887; "I" = 1.042893812997837; "J" = 1.045177571989127; "K" = 1.049140632416333; "L" = 1.046176491515212; "M" = 1.056159768562237; "N" = 1.046661844024895 }
    25 = @{ "B" = 1.02; "C" = 1.040652500712884; "D" = 1.047108107403013; "E" = 1.04438319325756; "F" = 1.054466258290173; "I" = 1.043385324375706; "J" = 1.046216248235232; "K" = 1.050125517502948; "L" = 1.047408990092219; "M" = 1.057461257631348; "N" = 1.047701995310407 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item([int]$rowNum, $colIndex[$col]).Value = $rowData[$col]
    }
}
